$p = $ppt.ActivePresentation

# --- Slide 52: "Computing Relative Addresses" -------------------------
# Move "we " from the start of run 2 to the end of run 1 in the first
# paragraph of the content placeholder:
#   "Similar to what " + "we did for "  ->  "Similar to what we " + "did for "
$s52 = $p.Slides.Item(52)
$sh52 = $s52.Shapes.Item(2)
$tr52 = $sh52.TextFrame.TextRange

$run1 = $tr52.Characters(1, 16)
$run1.Text = "Similar to what we "

$run2 = $tr52.Characters(20, 11)
$run2.Text = "did for "

# --- Slide 59: "Referencing Variables and Parameters for Function max()"
# Fix the opcode typo "LDDADDR -12" -> "LDLADDR -12" (loads the address of
# a local variable, not a "D" addr) while keeping "-12" as its own run.
$s59 = $p.Slides.Item(59)
$sh59 = $s59.Shapes.Item(2)
$tr59 = $sh59.TextFrame.TextRange

$opRun = $tr59.Characters(1, 8)
$opRun.Text = "LDLADDR "
